# gitPull.docx edit: rename the two byline paragraphs and fold the
# trailing (empty) bookmark-only paragraph into the second line.

$d = $word.ActiveDocument

# --- Paragraph 1: "Tyler " + (spell-check-flagged) "Roop" -> "Tyler Roop" ---
# Delete the whole paragraph (text + its paragraph mark) and retype it as a
# single run; this also clears the now-stale spellStart/spellEnd
# proofing-error bookmarks that wrapped "Roop".
$p1 = $d.Paragraphs(1)
$p1Start = $p1.Range.Start
$p1Full = $d.Range($p1Start, $p1.Range.End)
$p1Full.Delete()
$d.Range($p1Start, $p1Start).InsertBefore("Tyler Roop`r")

# --- Paragraph 2: "Kyle Jordan" -> "Weihan Huang" ---
[void]$d.Content.Find.Execute("Kyle Jordan", $true, $false, $false, $false, $false,
                               $true, 1, $false, "Weihan Huang", 2)

# --- Merge the now-empty 3rd paragraph (it only held the _GoBack bookmark)
# into paragraph 2 by removing the paragraph break between them. ---
$p2 = $d.Paragraphs(2)
$p2Mark = $d.Range($p2.Range.End - 1, $p2.Range.End)
$p2Mark.Delete()
